$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated feedback tallies in column C (new questionnaires collated in)
$values = @{
  2  = 0
  3  = 9
  4  = 10
  5  = 5
  6  = 9
  7  = 6
  9  = 0
  10 = 0
  11 = 1
  12 = 7
  13 = 2
  15 = 0
  16 = 0
  17 = 3
  18 = 6
  19 = 1
}

foreach ($r in $values.Keys) {
  $ws.Range("C$r").Value = $values[$r]
}

# Apply a thin box border around each label/value pair, and bold the labels.
# Rows 8 and 14 are intentionally blank separator rows and are skipped.
$rows = 2,3,4,5,6,7,9,10,11,12,13,15,16,17,18,19
foreach ($r in $rows) {
  $rowRange = $ws.Range("B$r`:C$r")
  $rowRange.Borders.LineStyle = 1
  $ws.Range("B$r").Font.Bold = $true
}

# Page setup for printing
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Move the active selection
$ws.Range("G17").Select() | Out-Null
